$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.465.64'
$ws.Range('E2').Value = '  +7.78%  '
$ws.Range('D3').Value = '3.626.26'
$ws.Range('E3').Value = '  +7.69%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = "'592.99"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.01%  '
$ws.Range('D6').Value = "'193.28"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +11.09%  '
$ws.Range('D7').Value = "'0.649"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.41%  '
$ws.Range('D8').Value = '3.600.64'
$ws.Range('E8').Value = '  +7.17%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').Value = "'0.182"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.51%  '
$ws.Range('D11').Value = "'0.664"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.43%  '
$ws.Range('D12').Value = "'58.04"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +9.87%  '
$ws.Range('D13').Value = "'0.0000297"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.93%  '
$ws.Range('D14').Value = "'9.78"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.11%  '
$ws.Range('D15').Value = '4.207.51'
$ws.Range('E15').Value = '  +7.90%  '
$ws.Range('D16').Value = '3.629.16'
$ws.Range('E16').Value = '  +7.72%  '
$ws.Range('E17').Value = '  +6.90%  '
$ws.Range('D18').Value = '70.297.61'
$ws.Range('E18').Value = '  +7.68%  '
$ws.Range('D19').Value = "'12.65"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.48%  '
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range('E21').Value = '  +6.34%  '
$ws.Range('D22').Value = "'496.64"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.40%  '
$ws.Range('E23').Value = '  +12.94%  '
$ws.Range('E24').Value = '  +18.67%  '
$ws.Range('E25').Value = '  +9.66%  '
$ws.Range('D26').Value = "'90.91"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').Value = "'3.12"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.26%  '
$ws.Range('D28').Value = "'11.24"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.00%  '
$ws.Range('D29').Value = "'9.42"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.13%  '
$ws.Range('D30').Value = "'32.40"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.52%  '
$ws.Range('D31').Value = "'7.56"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +16.52%  '
$ws.Range('D32').Value = "'12.25"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.85%  '
$ws.Range('D33').Value = "'617.59"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.49%  '
$ws.Range('D34').Value = "'65.43"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.62%  '
$ws.Range('E35').Value = '  +8.77%  '
$ws.Range('D36').Value = '0.0₃0836'
$ws.Range('E36').Value = '  +14.48%  '
$ws.Range('D37').Value = "'0.149"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.24%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = "'3.72"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.30%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = "'1.00"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').Value = "'38.04"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.30%  '
$ws.Range('E41').Value = '  +8.60%  '
$ws.Range('D42').Value = '3.338.69'
$ws.Range('E42').Value = '  +7.60%  '
$ws.Range('D43').Value = "'3.09"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +11.03%  '
$ws.Range('D44').Value = "'0.0448"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.37%  '
$ws.Range('E45').Value = '  +11.30%  '
$ws.Range('D46').Value = "'3.31"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.11%  '
$ws.Range('E47').Value = '  +3.66%  '
$ws.Range('D48').Value = "'2.78"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +16.65%  '
$ws.Range('D49').Value = "'9.09"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.91%  '
$ws.Range('D50').Value = "'3.27"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.57%  '
$ws.Range('E51').Value = '  -0.13%  '
